{"js": "// Ask only about user/other_party, not plaintiff/defendant:\n// rename the merge field `who_pays_new` -> `which_side_pays_new`\n// everywhere it appears in the document body.\nconst body = context.document.body;\nconst hits = body.search(\"who_pays_new\", { matchCase: true, matchWholeWord: false });\nhits.load(\"items\");\nawait context.sync();\n\nfor (let i = 0; i < hits.items.length; i++) {\n  hits.items[i].insertText(\"which_side_pays_new\", \"Replace\");\n}\nawait context.sync();\n", "ps1": "# Ask only about user/other_party, not plaintiff/defendant:\n# rename the merge field `who_pays_new` -> `which_side_pays_new`\n# everywhere it appears in the document body.\n$d = $word.ActiveDocument\n\n$find = $d.Content.Find\n$find.Text = \"who_pays_new\"\n$find.Replacement.Text = \"which_side_pays_new\"\n$find.Execute(\n    \"who_pays_new\",          # FindText\n    $false,                  # MatchCase\n    $false,                  # MatchWholeWord\n    $false,                  # MatchWildcards\n    $false,                  # MatchSoundsLike\n    $false,                  # MatchAllWordForms\n    $true,                   # Forward\n    1,                       # Wrap (wdFindContinue)\n    $false,                  # Format\n    \"which_side_pays_new\",   # ReplaceWith\n    2                        # Replace (wdReplaceAll)\n)\n"}
